$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H4").Value = 215
$ws.Range("I4").Value = 215
$ws.Range("K4").Value = 215
$ws.Range("M4").Value = -101
$ws.Range("H11").Value = 523.7646999999999
$ws.Range("I11").Value = 523.7646999999999
$ws.Range("K11").Value = 523.7646999999999
$ws.Range("M11").Value = -383.7646999999999
$ws.Range("H17").Value = 657028.6
$ws.Range("J17").Value = 657028.6
$ws.Range("L17").Value = 1971085.8
$ws.Range("N17").Value = -1971421.8
$ws.Range("H20").Value = 14624
$ws.Range("I20").Value = 2249.5
$ws.Range("J20").Value = 26998.5
$ws.Range("K20").Value = 2249.5
$ws.Range("L20").Value = 26998.5
$ws.Range("M20").Value = -2019.5
$ws.Range("N20").Value = -27458.5
$ws.Range("H35").Value = 14624
$ws.Range("I35").Value = 2249.5
$ws.Range("J35").Value = 26998.5
$ws.Range("K35").Value = 2249.5
$ws.Range("L35").Value = 26998.5
$ws.Range("M35").Value = -1870.5
$ws.Range("N35").Value = -27756.5
$ws.Range("H38").Value = 1277.75
$ws.Range("I38").Value = 1277.75
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3833.25
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -3461.25
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("H58").Value = 2700.8462
$ws.Range("J58").Value = 4917
$ws.Range("L58").Value = 14751
$ws.Range("N58").Value = -15051
$ws.Range("H112").Value = 108754.84
$ws.Range("J112").Value = 66249.875
$ws.Range("L112").Value = 198749.625
$ws.Range("N112").Value = -200965.625
$ws.Range("H113").Value = 2997
$ws.Range("I113").Value = 3139
$ws.Range("K113").Value = 3139
$ws.Range("M113").Value = 115
$ws.Range("H116").Value = 4692.846
$ws.Range("I116").Value = 4364.273
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 4364.273
$ws.Range("L116").Value = 6500
$ws.Range("M116").Value = -922.2730000000001
$ws.Range("N116").Value = -13384
$ws.Range("H121").Value = 2000
$ws.Range("J121").Value = 2000
$ws.Range("L121").Value = 6000
$ws.Range("N121").Value = -9494
$ws.Range("H125").Value = 2949.7144
$ws.Range("I125").Value = 2890.5454
$ws.Range("K125").Value = 26014.9086
$ws.Range("M125").Value = -23554.9086
$ws.Range("H131").Value = 2189
$ws.Range("I131").Value = 1147.1428
$ws.Range("K131").Value = 3441.4284
$ws.Range("M131").Value = 1598.5716
$ws.Range("H132").Value = 1190.579
$ws.Range("I132").Value = 1236.0555
$ws.Range("K132").Value = 3708.1665
$ws.Range("M132").Value = -1178.1665
$ws.Range("H137").Value = 2015.8445
$ws.Range("I137").Value = 1734.2122
$ws.Range("K137").Value = 5202.6366
$ws.Range("M137").Value = -2652.6366
$ws.Range("H138").Value = 2388.9429
$ws.Range("I138").Value = 1916.5652
$ws.Range("J138").Value = 3294.3333
$ws.Range("K138").Value = 5749.6956
$ws.Range("L138").Value = 9882.999899999999
$ws.Range("M138").Value = -609.6956
$ws.Range("N138").Value = -20162.9999
$ws.Range("H141").Value = 610.2105
$ws.Range("I141").Value = 585.9167
$ws.Range("J141").Value = 1047.5
$ws.Range("K141").Value = 1757.7501
$ws.Range("L141").Value = 3142.5
$ws.Range("M141").Value = 3422.2499
$ws.Range("N141").Value = -13502.5

foreach ($addr in @("N38", "N48", "N56")) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H2").Value = 1479077
$ws.Range("I2").Value = 1479077
$ws.Range("K2").Value = 1479077
$ws.Range("M2").Value = -1478964
$ws.Range("H4").Value = 1316.8334
$ws.Range("I4").Value = 1224.75
$ws.Range("K4").Value = 1224.75
$ws.Range("M4").Value = -1108.75
$ws.Range("H5").Value = 1205.2667
$ws.Range("I5").Value = 1203.5555
$ws.Range("K5").Value = 1203.5555
$ws.Range("M5").Value = -1091.5555
$ws.Range("H10").Value = 21797.285
$ws.Range("I10").Value = 19396.5
$ws.Range("J10").Value = 24998.334
$ws.Range("K10").Value = 19396.5
$ws.Range("L10").Value = 24998.334
$ws.Range("M10").Value = -19226.5
$ws.Range("N10").Value = -25338.334
$ws.Range("H32").Value = 3969.7026
$ws.Range("I32").Value = 4013.4722
$ws.Range("J32").Value = 2394
$ws.Range("K32").Value = 4013.4722
$ws.Range("L32").Value = 2394
$ws.Range("M32").Value = -3726.4722
$ws.Range("N32").Value = -2968
$ws.Range("H45").Value = 2159.6
$ws.Range("I45").Value = 2159.6
$ws.Range("K45").Value = 2159.6
$ws.Range("M45").Value = -1782.6
$ws.Range("H61").Value = 76925160
$ws.Range("I61").Value = 100001760
$ws.Range("J61").Value = 3165.6667
$ws.Range("K61").Value = 100001760
$ws.Range("L61").Value = 3165.6667
$ws.Range("M61").Value = -100001548
$ws.Range("N61").Value = -3589.6667
$ws.Range("H74").Value = 43487976
$ws.Range("I74").Value = 66676732
$ws.Range("K74").Value = 66676732
$ws.Range("M74").Value = -66675858
$ws.Range("H77").Value = 43487976
$ws.Range("I77").Value = 66676732
$ws.Range("K77").Value = 333383660
$ws.Range("M77").Value = -333379292
$ws.Range("H96").Value = 46637
$ws.Range("J96").Value = 46637
$ws.Range("L96").Value = 46637
$ws.Range("N96").Value = -52129
$ws.Range("H97").Value = 368.25
$ws.Range("I97").Value = 392.2857
$ws.Range("K97").Value = 392.2857
$ws.Range("M97").Value = 103.7143
$ws.Range("H102").Value = 12503854
$ws.Range("I102").Value = 14289446
$ws.Range("J102").Value = 4709
$ws.Range("K102").Value = 14289446
$ws.Range("L102").Value = 4709
$ws.Range("M102").Value = -14287824
$ws.Range("N102").Value = -7953
$ws.Range("H116").Value = 1479077
$ws.Range("I116").Value = 1479077
$ws.Range("K116").Value = 1479077
$ws.Range("M116").Value = -1476783
$ws.Range("H122").Value = 3154.8
$ws.Range("I122").Value = 3118.5
$ws.Range("J122").Value = 3300
$ws.Range("K122").Value = 9355.5
$ws.Range("L122").Value = 9900
$ws.Range("M122").Value = -6905.5
$ws.Range("N122").Value = -14800
$ws.Range("H132").Value = 3849006
$ws.Range("I132").Value = 3849006
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11547018
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -11544488
$ws.Range("H136").Value = 76925160
$ws.Range("I136").Value = 100001760
$ws.Range("J136").Value = 3165.6667
$ws.Range("K136").Value = 300005280
$ws.Range("L136").Value = 9497.000100000001
$ws.Range("M136").Value = -300002730
$ws.Range("N136").Value = -14597.0001

foreach ($addr in @("N132")) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H3").Value = 1479077
$ws.Range("I3").Value = 1479077
$ws.Range("K3").Value = 1479077
$ws.Range("M3").Value = -1478963
$ws.Range("H4").Value = 1205.2667
$ws.Range("I4").Value = 1203.5555
$ws.Range("K4").Value = 1203.5555
$ws.Range("M4").Value = -1088.5555
$ws.Range("H29").Value = 13000
$ws.Range("I29").Value = 13000
$ws.Range("K29").Value = 13000
$ws.Range("M29").Value = -12711
$ws.Range("H99").Value = 1579.85
$ws.Range("I99").Value = 1270.2858
$ws.Range("K99").Value = 1270.2858
$ws.Range("M99").Value = 227.7141999999999
$ws.Range("H105").Value = 2451.2
$ws.Range("I105").Value = 2497.5
$ws.Range("J105").Value = 2266
$ws.Range("K105").Value = 2497.5
$ws.Range("L105").Value = 2266
$ws.Range("M105").Value = -750.5
$ws.Range("N105").Value = -5760
$ws.Range("H134").Value = 34335084
$ws.Range("I134").Value = 39616972
$ws.Range("K134").Value = 118850916
$ws.Range("M134").Value = -118848381

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 838168.4
$ws.Range("I16").Value = 1087830.4
$ws.Range("J16").Value = 5961.6665
$ws.Range("K16").Value = 1087830.4
$ws.Range("L16").Value = 5961.6665
$ws.Range("M16").Value = -1087543.4
$ws.Range("N16").Value = -6535.6665
$ws.Range("H31").Value = 6174.521
$ws.Range("I31").Value = 3558.1843
$ws.Range("J31").Value = 16116.6
$ws.Range("K31").Value = 3558.1843
$ws.Range("L31").Value = 16116.6
$ws.Range("M31").Value = -3263.1843
$ws.Range("N31").Value = -16706.6
$ws.Range("H34").Value = 6174.521
$ws.Range("I34").Value = 3558.1843
$ws.Range("J34").Value = 16116.6
$ws.Range("K34").Value = 3558.1843
$ws.Range("L34").Value = 16116.6
$ws.Range("M34").Value = -3356.1843
$ws.Range("N34").Value = -16520.6
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("H58").Value = 20004820
$ws.Range("I58").Value = 26321684
$ws.Range("J58").Value = 1416.1666
$ws.Range("K58").Value = 26321684
$ws.Range("L58").Value = 1416.1666
$ws.Range("M58").Value = -26321481
$ws.Range("N58").Value = -1822.1666
$ws.Range("H97").Value = 39998.5
$ws.Range("J97").Value = 39998.5
$ws.Range("L97").Value = 39998.5
$ws.Range("N97").Value = -41980.5
$ws.Range("H98").Value = 107499.5
$ws.Range("J98").Value = 107499.5
$ws.Range("L98").Value = 107499.5
$ws.Range("N98").Value = -111991.5
$ws.Range("H100").Value = 48000
$ws.Range("J100").Value = 48000
$ws.Range("L100").Value = 48000
$ws.Range("N100").Value = -50164
$ws.Range("H105").Value = 1112346.9
$ws.Range("I105").Value = 1539359.8
$ws.Range("K105").Value = 1539359.8
$ws.Range("M105").Value = -1537612.8
$ws.Range("H113").Value = 838168.4
$ws.Range("I113").Value = 1087830.4
$ws.Range("J113").Value = 5961.6665
$ws.Range("K113").Value = 1087830.4
$ws.Range("L113").Value = 5961.6665
$ws.Range("M113").Value = -1085660.4
$ws.Range("N113").Value = -10301.6665
$ws.Range("H132").Value = 25001238
$ws.Range("J132").Value = 1023.4
$ws.Range("L132").Value = 3070.2
$ws.Range("N132").Value = -8130.2
$ws.Range("H134").Value = 8622645
$ws.Range("I134").Value = 11365664
$ws.Range("J134").Value = 1728
$ws.Range("K134").Value = 34096992
$ws.Range("L134").Value = 5184
$ws.Range("M134").Value = -34094457
$ws.Range("N134").Value = -10254
$ws.Range("H136").Value = 20004820
$ws.Range("I136").Value = 26321684
$ws.Range("J136").Value = 1416.1666
$ws.Range("K136").Value = 78965052
$ws.Range("L136").Value = 4248.4998
$ws.Range("M136").Value = -78962502
$ws.Range("N136").Value = -9348.4998

foreach ($addr in @("M55")) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3346
$ws.Range("H36").Value = 523.5
$ws.Range("I36").Value = 531.3333
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 1593.9999
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = -1424.9999
$ws.Range("N36").Value = -1838
$ws.Range("H37").Value = 149751
$ws.Range("J37").Value = 149751
$ws.Range("L37").Value = 449253
$ws.Range("N37").Value = -449477
$ws.Range("H38").Value = 180.42857
$ws.Range("J38").Value = 118.2
$ws.Range("L38").Value = 354.6
$ws.Range("N38").Value = -1048.6
$ws.Range("H49").Value = 7999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 7999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 23997
$ws.Range("N49").Value = -24309
$ws.Range("H50").Value = 3883.5
$ws.Range("I50").Value = 3326.25
$ws.Range("K50").Value = 9978.75
$ws.Range("M50").Value = -9497.75
$ws.Range("H53").Value = 3883.5
$ws.Range("I53").Value = 3326.25
$ws.Range("K53").Value = 9978.75
$ws.Range("M53").Value = -9497.75
$ws.Range("H55").Value = 99
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 99
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 297
$ws.Range("N55").Value = -651
$ws.Range("H75").Value = 2000437.4
$ws.Range("J75").Value = 10000000
$ws.Range("L75").Value = 30000000
$ws.Range("N75").Value = -30001996
$ws.Range("H78").Value = 2000437.4
$ws.Range("J78").Value = 10000000
$ws.Range("L78").Value = 90000000
$ws.Range("N78").Value = -90009984
$ws.Range("H92").Value = 366.33334
$ws.Range("I92").Value = 366.33334
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1099.00002
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 148.9999800000001
$ws.Range("H129").Value = 2393.6155
$ws.Range("J129").Value = 4299.6665
$ws.Range("L129").Value = 12898.9995
$ws.Range("N129").Value = -22898.9995
$ws.Range("H132").Value = 3622.8
$ws.Range("I132").Value = 2034.6666
$ws.Range("K132").Value = 18311.9994
$ws.Range("M132").Value = -15781.9994

foreach ($addr in @("M16", "M49", "M55", "N92")) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H19").Value = 2252.5
$ws.Range("I19").Value = 2252.5
$ws.Range("K19").Value = 2252.5
$ws.Range("M19").Value = -1964.5
$ws.Range("H63").Value = 39000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H64").Value = 59950
$ws.Range("J64").Value = 59950
$ws.Range("L64").Value = 59950
$ws.Range("N64").Value = -60446
$ws.Range("H66").Value = 39000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H67").Value = 59950
$ws.Range("J67").Value = 59950
$ws.Range("L67").Value = 59950
$ws.Range("N67").Value = -61666
$ws.Range("H70").Value = 4191.125
$ws.Range("I70").Value = 4085.1667
$ws.Range("K70").Value = 4085.1667
$ws.Range("M70").Value = -3815.1667
$ws.Range("H73").Value = 4191.125
$ws.Range("I73").Value = 4085.1667
$ws.Range("K73").Value = 4085.1667
$ws.Range("M73").Value = -3149.1667
$ws.Range("H80").Value = 4496
$ws.Range("I80").Value = 3995
$ws.Range("K80").Value = 3995
$ws.Range("M80").Value = -2997
$ws.Range("H83").Value = 4496
$ws.Range("I83").Value = 3995
$ws.Range("K83").Value = 19975
$ws.Range("M83").Value = -14983
$ws.Range("H98").Value = 117600
$ws.Range("J98").Value = 117600
$ws.Range("L98").Value = 117600
$ws.Range("N98").Value = -123590
$ws.Range("H99").Value = 17831.5
$ws.Range("I99").Value = 7398
$ws.Range("K99").Value = 7398
$ws.Range("M99").Value = -5152
$ws.Range("H100").Value = 119999
$ws.Range("J100").Value = 119999
$ws.Range("L100").Value = 119999
$ws.Range("N100").Value = -122163
$ws.Range("H107").Value = 1714.3684
$ws.Range("I107").Value = 1531.8334
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1531.8334
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = 388.1666
$ws.Range("N107").Value = -8840
$ws.Range("H108").Value = 89999.5
$ws.Range("J108").Value = 89999.5
$ws.Range("L108").Value = 89999.5
$ws.Range("N108").Value = -97679.5
$ws.Range("H122").Value = 6908.923
$ws.Range("I122").Value = 3706.8572
$ws.Range("K122").Value = 11120.5716
$ws.Range("M122").Value = -8670.571599999999
$ws.Range("H126").Value = 4309.25
$ws.Range("I126").Value = 4629.0713
$ws.Range("K126").Value = 13887.2139
$ws.Range("M126").Value = -11417.2139
$ws.Range("H132").Value = 7816576
$ws.Range("I132").Value = 9618858
$ws.Range("K132").Value = 28856574
$ws.Range("M132").Value = -28854044

foreach ($addr in @("N63", "N66")) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H7").Value = 3788.8096
$ws.Range("I7").Value = 3808.7334
$ws.Range("J7").Value = 3739
$ws.Range("K7").Value = 3808.7334
$ws.Range("L7").Value = 3739
$ws.Range("M7").Value = -3696.7334
$ws.Range("N7").Value = -3963
$ws.Range("H12").Value = 2799.6667
$ws.Range("I12").Value = 2699.5
$ws.Range("J12").Value = 3000
$ws.Range("K12").Value = 2699.5
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = -2529.5
$ws.Range("N12").Value = -3340
$ws.Range("H61").Value = 2587.0588
$ws.Range("J61").Value = 2675.4
$ws.Range("L61").Value = 2675.4
$ws.Range("N61").Value = -3079.4
$ws.Range("H103").Value = 25866.334
$ws.Range("J103").Value = 25866.334
$ws.Range("L103").Value = 25866.334
$ws.Range("N103").Value = -28210.334
$ws.Range("H113").Value = 2587.0588
$ws.Range("J113").Value = 2675.4
$ws.Range("L113").Value = 2675.4
$ws.Range("N113").Value = -7015.4
$ws.Range("H126").Value = 3788.8096
$ws.Range("I126").Value = 3808.7334
$ws.Range("J126").Value = 3739
$ws.Range("K126").Value = 11426.2002
$ws.Range("L126").Value = 11217
$ws.Range("M126").Value = -8956.200199999999
$ws.Range("N126").Value = -16157
$ws.Range("H132").Value = 34296132
$ws.Range("I132").Value = 43648560
$ws.Range("K132").Value = 130945680
$ws.Range("M132").Value = -130943150
$ws.Range("H136").Value = 2434.5625
$ws.Range("J136").Value = 2719.923
$ws.Range("L136").Value = 8159.768999999999
$ws.Range("N136").Value = -13259.769

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("H74").Value = 13169.875
$ws.Range("I74").Value = 15997.5
$ws.Range("J74").Value = 12227.333
$ws.Range("K74").Value = 15997.5
$ws.Range("L74").Value = 12227.333
$ws.Range("M74").Value = -15061.5
$ws.Range("N74").Value = -14099.333
$ws.Range("H75").Value = 69999
$ws.Range("J75").Value = 69999
$ws.Range("L75").Value = 69999
$ws.Range("N75").Value = -71871
$ws.Range("H77").Value = 13169.875
$ws.Range("I77").Value = 15997.5
$ws.Range("J77").Value = 12227.333
$ws.Range("K77").Value = 47992.5
$ws.Range("L77").Value = 36681.999
$ws.Range("M77").Value = -43312.5
$ws.Range("N77").Value = -46041.999
$ws.Range("H78").Value = 69999
$ws.Range("J78").Value = 69999
$ws.Range("L78").Value = 209997
$ws.Range("N78").Value = -219357
$ws.Range("H86").Value = 25000
$ws.Range("J86").Value = 25000
$ws.Range("L86").Value = 25000
$ws.Range("N86").Value = -27246
$ws.Range("H89").Value = 25000
$ws.Range("J89").Value = 25000
$ws.Range("L89").Value = 125000
$ws.Range("N89").Value = -136232
$ws.Range("H132").Value = 8932218
$ws.Range("I132").Value = 10639994
$ws.Range("K132").Value = 31919982
$ws.Range("M132").Value = -31917452
$ws.Range("H136").Value = 9617727
$ws.Range("I136").Value = 10419046
$ws.Range("J136").Value = 1899
$ws.Range("K136").Value = 31257138
$ws.Range("M136").Value = -31254588
$ws.Range("N136").Value = -10797

foreach ($addr in @("N51")) {
    $ws.Range($addr).ClearContents()
}
